# Update the simulated "realeffort" values (column F) for every worker row,
# and swap the prolificid/name/gender for the two workers whose realeffort
# values crossed over (rows 9 and 10), matching the new ranking data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New realeffort (column F) values for rows 2-13 ---
$ws.Cells.Item(2, 6).Value  = 11.12356095231806
$ws.Cells.Item(3, 6).Value  = 10.07657103797102
$ws.Cells.Item(4, 6).Value  = 8.469824362969149
$ws.Cells.Item(5, 6).Value  = 8.252431536799262
$ws.Cells.Item(6, 6).Value  = 7.094657342882389
$ws.Cells.Item(7, 6).Value  = 6.351992923050718
$ws.Cells.Item(8, 6).Value  = 6.111943368614604
$ws.Cells.Item(9, 6).Value  = 5.465857846036377
$ws.Cells.Item(10, 6).Value = 5.411049145544538
$ws.Cells.Item(11, 6).Value = 4.011565163053068
$ws.Cells.Item(12, 6).Value = 2.037811163075423
$ws.Cells.Item(13, 6).Value = 0.2234880371687213

# --- Rows 9 and 10 swap prolificid (C), name (D) and gender (E), ---
# --- and their "index" source id (column B) swaps along with them. ---
$ws.Cells.Item(9, 2).Value  = 3
$ws.Cells.Item(9, 3).Value  = "60bd88b8fc436774352f53b9"
$ws.Cells.Item(9, 4).Value  = "Annes"
$ws.Cells.Item(9, 5).Value  = "female"

$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = "5c27de12a2b00a00018b2c16"
$ws.Cells.Item(10, 4).Value = "Ankai"
$ws.Cells.Item(10, 5).Value = "male"
